$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking values (e.g. "0.9994")
# are not auto-converted to floating point numbers, matching original inline-string content.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "22.081.01"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "1.558.32"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "0.9997"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "287.81"
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("D7").Value = "0.3878"
$ws.Range("E7").Value = "  +4.46%  "
$ws.Range("D8").Value = "0.3246"
$ws.Range("E8").Value = "  -0.99%  "
$ws.Range("D9").Value = "42.90"
$ws.Range("E9").Value = "  -7.73%  "
$ws.Range("D10").Value = "1.126"
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("D11").Value = "0.07378"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").Value = "0.9995"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "19.39"
$ws.Range("E13").Value = "  -5.22%  "
$ws.Range("D14").Value = "5.708"
$ws.Range("E14").Value = "  -2.33%  "
$ws.Range("D15").Value = "6.811"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").Value = "0.00001132"
$ws.Range("E16").Value = "  +2.93%  "
$ws.Range("D17").Value = "1.558.10"
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("D18").Value = "0.06612"
$ws.Range("E18").Value = "  -1.30%  "
$ws.Range("D19").Value = "85.40"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").Value = "6.404"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").Value = "0.9991"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "16.01"
$ws.Range("E22").Value = "  -1.53%  "
$ws.Range("D23").Value = "11.50"
$ws.Range("E23").Value = "  -2.35%  "
$ws.Range("D24").Value = "22.079.46"
$ws.Range("E24").Value = "  -1.42%  "
$ws.Range("D25").Value = "2.341"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("D26").Value = "2.562"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").Value = "149.63"
$ws.Range("E27").Value = "  -0.80%  "
$ws.Range("D28").Value = "18.91"
$ws.Range("E28").Value = "  -2.28%  "
$ws.Range("D29").Value = "4.868"
$ws.Range("E29").Value = "  -1.63%  "
$ws.Range("D30").Value = "1.728.77"
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("D31").Value = "121.01"
$ws.Range("E31").Value = "  -2.13%  "
$ws.Range("D32").Value = "1.113"
$ws.Range("E32").Value = "  +6.00%  "
$ws.Range("D33").Value = "5.859"
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("D34").Value = "1.709"
$ws.Range("E34").Value = "  -13.00%  "
$ws.Range("D35").Value = "0.08208"
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("D36").Value = "9.298"
$ws.Range("E36").Value = "  -3.94%  "
$ws.Range("D37").Value = "0.06272"
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("D38").Value = "0.02309"
$ws.Range("E38").Value = "  -3.13%  "
$ws.Range("D39").Value = "5.236"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "0.2113"
$ws.Range("E40").Value = "  -3.39%  "
$ws.Range("D41").Value = "1.223"
$ws.Range("E41").Value = "  -6.95%  "
$ws.Range("D42").Value = "10.90"
$ws.Range("E42").Value = "  -1.96%  "
$ws.Range("D43").Value = "0.9989"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").Value = "0.5969"
$ws.Range("E44").Value = "  -2.39%  "
$ws.Range("D45").Value = "13.51"
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("D46").Value = "3.718"
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("D47").Value = "0.5765"
$ws.Range("E47").Value = "  -3.10%  "
$ws.Range("D48").Value = "1.932"
$ws.Range("E48").Value = "  -3.88%  "
$ws.Range("D49").Value = "119.18"
$ws.Range("E49").Value = "  -3.66%  "
$ws.Range("D50").Value = "1.161"
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("D51").Value = "0.06904"
$ws.Range("E51").Value = "  -3.47%  "
